$wb = $excel.ActiveWorkbook

# --- Sheet "Alluvial for Mapping": fix "Mortendad" -> "Mortandad" typo in column P, rows 6-21 ---
$wsMap = $wb.Worksheets.Item("Alluvial for Mapping")
for ($r = 6; $r -le 21; $r++) {
    $wsMap.Cells.Item($r, 16).Value = "Mortandad"
}

# --- Sheet "Alluvial Exhibit": correct watershed naming / headings ---
$wsExh = $wb.Worksheets.Item("Alluvial Exhibit")

# Fix "Mortendad Canyon" -> "Mortandad Canyon" heading
$wsExh.Range("A8").Value = "Mortandad Canyon"

# Split "Los Alamos and Pajarito Canyons" into its own "Los Alamos Canyon" section
$wsExh.Range("A25").Value = "Los Alamos Canyon"

# Insert a new header row above the Pajarito Canyon wells (currently starting at row 44),
# and merge it across A:H like the other section-header rows (merge first, so the later
# per-cell style copy below isn't re-split across the merged span)
$wsExh.Rows.Item(44).Insert()
$wsExh.Range("A44:H44").Merge()

# Copy the formatting of an existing section header row (row 8: "Mortandad Canyon") onto the
# newly inserted row (A = section-title style, H = matching right-edge style), clear the
# leftover formatting Excel auto-applied to the inserted row's middle cells, then set the
# new header's text
$wsExh.Range("A8").Copy()
$wsExh.Range("A44").PasteSpecial(-4122)
$wsExh.Range("H8").Copy()
$wsExh.Range("H44").PasteSpecial(-4122)
$wsExh.Range("B44:G44").ClearFormats()
$wsExh.Range("A44").Value = "Pajarito Canyon"

# Widen column C slightly (stored width 12 -> 13)
$wsExh.Columns.Item(3).ColumnWidth = 12.2

Write-Host "edit complete"
